$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Header / title font restyle -----------------------------------------
# Title (A1) keeps Bold but drops the explicit 14pt size (back to default)
# and becomes white, matching the new shared "bold white" font.
$ws1.Range("A1").Font.Size = 11
$ws1.Range("A1").Font.Bold = $true
$ws1.Range("A1").Font.Color = 16777215

# Header row (A2:K2) keeps its dark-blue fill/border but the bold font
# becomes white too, so it now shares the very same font as the title.
$ws1.Range("A2:K2").Font.Bold = $true
$ws1.Range("A2:K2").Font.Color = 16777215

# --- Data updates ----------------------------------------------------------
$ws1.Range("H3").Value = 518
$ws1.Range("I3").Value = "'16-Sep-2025"

$ws1.Range("H4").Value = -55
$ws1.Range("I4").Value = "'16-Sep-2025"
